$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing data block (rows 2-21) down by 4 rows, to make
# room for 4 new rows under the header (avoids Range.Insert(), which
# would otherwise drag the header's bold/centered style onto the shifted
# cells). Walk bottom-up so we never clobber a row before reading it.
for ($r = 21; $r -ge 2; $r--) {
    $dest = $r + 4
    $ws.Cells.Item($dest, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

# --- Fill the 4 newly freed rows (2-5) with the new data.
$topRows = @(
    @(-0.0310014113783836, 0.0068722339347004, 0.001527163083665),
    @(-0.0077885319478809, -0.0080939643085002, -0.0400116741657257),
    @(-0.0332921557128429, -0.0244346093386411, -0.0705549344420433),
    @(0.0039706239476799, 0.0195476878434419, -0.038026362657547)
)

for ($i = 0; $i -lt $topRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value2 = $topRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $topRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $topRows[$i][2]
}

# --- Append 6 new rows at the bottom (rows 26-31).
$bottomRows = @(
    @(0.1531744599342346, 0.9390525817871094, -0.3572034537792206),
    @(0.1224784851074218, -0.0070249503478407, -0.1852448880672454),
    @(-0.0273362193256616, -0.1533271819353103, -0.1511891484260559),
    @(-0.0138971842825412, -0.0740674138069152, -0.052381694316864),
    @(-0.102472648024559, -0.0630718395113945, -0.0003054326225537),
    @(0.0174096599221229, 0.0529925599694252, -0.0123700210824608)
)

for ($i = 0; $i -lt $bottomRows.Count; $i++) {
    $r = 26 + $i
    $ws.Cells.Item($r, 1).Value2 = $bottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $bottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $bottomRows[$i][2]
}
